# Setup.pptx slide 2 ("Setup") content updates:
#  - "Download SPARK + SPARK_HOME + PATH" -> "Download SPARK 2.4 + SPARK_HOME + PATH"
#  - " IDEA + Scala plugin (File -> Settings " + "-> Plugins)" merged/updated to
#    " IDEA 2019 + Scala plugin (File -> Settings -> Plugins)" (single run)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# --- 1) "Download SPARK ..." paragraph (3rd paragraph) ---
$downloadPara = $tr.Paragraphs(3)
$downloadRun = $downloadPara.Runs(1)
$downloadRun.Text = "Download SPARK 2.4 + SPARK_HOME + PATH"

# --- 2) "Intellij IDEA ... -> Plugins)" paragraph (5th / last paragraph) ---
$ideaPara = $tr.Paragraphs(5)
$ideaRun2 = $ideaPara.Runs(2)
$ideaRun3 = $ideaPara.Runs(3)

# Merge run2 + run3 into a single run carrying the new text.
$combinedStart = $ideaRun2.Start
$combinedLength = ($ideaPara.Start + $ideaPara.Length) - $combinedStart
$combinedRange = $tr.Characters($combinedStart, $combinedLength)
$combinedRange.Text = " IDEA 2019 + Scala plugin (File -> Settings -> Plugins)"
